$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = -7.491600000000002
$ws.Range("C3").Value = -10.45179999999999
$ws.Range("D3").Value = -6.60459999999999
$ws.Range("C4").Value = -12.31199999999999
$ws.Range("E8").Value = 16.6802
$ws.Range("D9").Value = -7.207999999999999
$ws.Range("B11").Value = 5.936599999999998
$ws.Range("E11").Value = 16.29409999999999
$ws.Range("B12").Value = 4.793199999999996
$ws.Range("C14").Value = -13.91039999999999
$ws.Range("E14").Value = 16.8902
$ws.Range("B15").Value = 4.918099999999998
$ws.Range("D15").Value = -8.758299999999993
$ws.Range("E15").Value = 16.5498
$ws.Range("E17").Value = 16.83100000000001
$ws.Range("D19").Value = -8.104999999999997
$ws.Range("D20").Value = -7.858900000000001
$ws.Range("D25").Value = -7.775200000000009
$ws.Range("C26").Value = -12.0892
$ws.Range("E26").Value = 15.4825
$ws.Range("B27").Value = 6.3607
$ws.Range("D27").Value = -9.086099999999998
$ws.Range("B28").Value = 6.403000000000002
$ws.Range("D28").Value = -8.026800000000001
$ws.Range("D30").Value = -7.230000000000005
$ws.Range("B31").Value = 3.998199999999998
$ws.Range("C31").Value = -13.57149999999999
$ws.Range("B32").Value = 6.535199999999999
$ws.Range("D32").Value = -8.531500000000001
$ws.Range("C35").Value = -12.3351
$ws.Range("B36").Value = 9.268699999999997
$ws.Range("E36").Value = 16.3112
$ws.Range("C37").Value = -13.30479999999999
$ws.Range("B38").Value = 5.5542
$ws.Range("C39").Value = -12.2309
$ws.Range("C40").Value = -13.9242
$ws.Range("E42").Value = 16.48349999999999
$ws.Range("D44").Value = -7.275100000000006
$ws.Range("C45").Value = -13.7461
$ws.Range("B46").Value = 7.500100000000003
$ws.Range("D47").Value = -7.432
$ws.Range("C52").Value = -10.9223
$ws.Range("B54").Value = 5.003200000000003
$ws.Range("B55").Value = 5.073399999999997
$ws.Range("B56").Value = 4.735000000000002
$ws.Range("C57").Value = -14.40779999999999
$ws.Range("D58").Value = -7.791899999999998
$ws.Range("D62").Value = -8.488699999999993
$ws.Range("E64").Value = 17.4539
$ws.Range("B67").Value = 4.929099999999994
$ws.Range("E68").Value = 16.76360000000001
$ws.Range("B69").Value = 5.699199999999993
$ws.Range("B72").Value = 5.698900000000001
$ws.Range("B73").Value = 8.334899999999996
$ws.Range("D77").Value = -5.652299999999999
$ws.Range("D78").Value = -7.760600000000003
$ws.Range("E79").Value = 18.28680000000002
$ws.Range("C81").Value = -13.5194
$ws.Range("B83").Value = 5.520299999999997
$ws.Range("C83").Value = -13.83599999999999
$ws.Range("D84").Value = -8.7235
$ws.Range("B86").Value = 4.882800000000004
$ws.Range("D89").Value = -6.115899999999999
$ws.Range("E89").Value = 19.10900000000003
$ws.Range("B91").Value = 5.0771
$ws.Range("D91").Value = -6.123699999999999
$ws.Range("D92").Value = -6.041099999999999
$ws.Range("B93").Value = 6.872600000000003
$ws.Range("D96").Value = -7.801600000000013
$ws.Range("B99").Value = 4.658799999999997
$ws.Range("C100").Value = -12.8842
$ws.Range("C102").Value = -13.2561
$ws.Range("D102").Value = -7.556999999999996
